# Generate Report for Handoff
# Adds two new source files (and their per-language handoff records) to the
# localization-status report: one row each on the "Overview", "zh-cn" and
# "de-de" sheets, for files:
#   7472a77c-e5fd-479e-a3a9-956a736b8f16.md
#   ab9e920e-1d34-45e8-a965-90abe4b25f2b.md

$wb = $excel.ActiveWorkbook

$file3 = "7472a77c-e5fd-479e-a3a9-956a736b8f16"
$file4 = "ab9e920e-1d34-45e8-a965-90abe4b25f2b"

$xlf3 = "7472a77c-e5fd-479e-a3a9-956a736b8f16.659df0150749e6ff0145a1be6fe74fcec5e5be17"
$xlf4 = "ab9e920e-1d34-45e8-a965-90abe4b25f2b.0bcdeed460db9df57d6b49ffc783a0164cf80855"

$status = "Ready for handoff"
$handoffDateZh = "2016-03-20 16:45:33"
$handoffDateDe = "2016-03-20 16:45:41"
$noHandback = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(4, 2).Value = $status
$wsOverview.Cells.Item(4, 3).Value = $status
$wsOverview.Cells.Item(4, 4).Value = $handoffDateDe
$wsOverview.Cells.Item(4, 4).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Hyperlinks.Add(
    $wsOverview.Cells.Item(4, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/7f59539f4fdb1d399b164bc0e1a58231be1c75e8/e2e/$file3.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "$file3.md"
) | Out-Null

$wsOverview.Cells.Item(5, 2).Value = $status
$wsOverview.Cells.Item(5, 3).Value = $status
$wsOverview.Cells.Item(5, 4).Value = $handoffDateDe
$wsOverview.Cells.Item(5, 4).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Hyperlinks.Add(
    $wsOverview.Cells.Item(5, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/7f59539f4fdb1d399b164bc0e1a58231be1c75e8/e2e/$file4.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "$file4.md"
) | Out-Null

# ---------------------------------------------------------------------
# Sheets "zh-cn" / "de-de":
#  A Source File Name | B File Extension | C Status | D Latest Handoff File
#  E Latest Handoff Datetime | F Latest Target File | G Latest Handback File
#  H Latest Handback DateTime | I Reference Tokens | J Handoff Reason
#  K Dependency From | L Error Detail
# ---------------------------------------------------------------------
function Add-LangRows($ws, $lang, $handoffDate, $urlFilesBase, $urlXlfBase) {
    $row = 4
    foreach ($pair in @(@($file3, $xlf3), @($file4, $xlf4))) {
        $src = $pair[0]
        $xlf = $pair[1]

        $ws.Cells.Item($row, 2).Value = ".md"
        $ws.Cells.Item($row, 3).Value = $status
        $ws.Cells.Item($row, 5).Value = $handoffDate
        $ws.Cells.Item($row, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
        $ws.Cells.Item($row, 8).Value = $noHandback
        $ws.Cells.Item($row, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
        $ws.Cells.Item($row, 10).Value = "Include"

        $ws.Hyperlinks.Add(
            $ws.Cells.Item($row, 1),
            "$urlFilesBase/$src.md",
            [System.Reflection.Missing]::Value,
            [System.Reflection.Missing]::Value,
            "$src.md"
        ) | Out-Null

        $ws.Hyperlinks.Add(
            $ws.Cells.Item($row, 4),
            "$urlXlfBase/$xlf.$lang.xlf",
            [System.Reflection.Missing]::Value,
            [System.Reflection.Missing]::Value,
            "$xlf.$lang.xlf"
        ) | Out-Null

        $row++
    }
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$zhFilesBase = "https://github.com/OpenLocalizationTest/oltest/blob/7f59539f4fdb1d399b164bc0e1a58231be1c75e8/e2e"
$zhXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/06404e5c67abd08e398bec2fea30c42655f13285/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht"
Add-LangRows $wsZhCn "zh-cn" $handoffDateZh $zhFilesBase $zhXlfBase

$wsDeDe = $wb.Worksheets.Item("de-de")
$deFilesBase = "https://github.com/OpenLocalizationTest/oltest/blob/7f59539f4fdb1d399b164bc0e1a58231be1c75e8/e2e"
$deXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ab2bd34ba3d0b8884cf8fdd7fe73aca769fbbd04/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht"
Add-LangRows $wsDeDe "de-de" $handoffDateDe $deFilesBase $deXlfBase
